$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet shrinks from 23 rows to 21: drop the last two rows entirely
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()

# Rewrite rows 1-21 so each cell and row height matches the new layout
$ws.Cells.Item(1, 1).Value = $null
$ws.Cells.Item(1, 2).Value = "Ementa atual:"
$ws.Cells.Item(1, 3).Value = "Ementa modificada (dados modificados em vermelho):"
$ws.Rows.Item(1).RowHeight = 15

$ws.Cells.Item(2, 1).Value = $null
$ws.Cells.Item(2, 2).Value = "LOM3248"
$ws.Cells.Item(2, 3).Value = "LOM3248"
$ws.Rows.Item(2).RowHeight = 15

$ws.Cells.Item(3, 1).Value = "Nome:"
$ws.Cells.Item(3, 2).Value = " Tópicos Especiais em Engenharia Física I"
$ws.Cells.Item(3, 3).Value = " Tópicos Especiais em Engenharia Física I"
$ws.Rows.Item(3).RowHeight = 15

$ws.Cells.Item(4, 1).Value = "Name:"
$ws.Cells.Item(4, 2).Value = "Special Topics in Engineering Physics I"
$ws.Cells.Item(4, 3).Value = "Special Topics in Engineering Physics I"
$ws.Rows.Item(4).RowHeight = 15

$ws.Cells.Item(5, 1).Value = "Créditos-aula:"
$ws.Cells.Item(5, 2).Value = "4"
$ws.Cells.Item(5, 3).Value = "4"
$ws.Rows.Item(5).RowHeight = 15

$ws.Cells.Item(6, 1).Value = "Créditos-trabalho"
$ws.Cells.Item(6, 2).Value = "0"
$ws.Cells.Item(6, 3).Value = "0"
$ws.Rows.Item(6).RowHeight = 15

$ws.Cells.Item(7, 1).Value = "Carga horária:"
$ws.Cells.Item(7, 2).Value = "60 h"
$ws.Cells.Item(7, 3).Value = "60 h"
$ws.Rows.Item(7).RowHeight = 15

$ws.Cells.Item(8, 1).Value = "Ativação:"
$ws.Cells.Item(8, 2).Value = "01/01/2012"
$ws.Cells.Item(8, 3).Value = "01/01/2012"
$ws.Rows.Item(8).RowHeight = 15

$ws.Cells.Item(9, 1).Value = "Semestre ideal:"
$ws.Cells.Item(9, 2).Value = "EF-7"
$ws.Cells.Item(9, 3).Value = "EF-7"
$ws.Rows.Item(9).RowHeight = 15

$ws.Cells.Item(10, 1).Value = "Objetivos:"
$ws.Cells.Item(10, 2).Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Cells.Item(10, 3).Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Rows.Item(10).RowHeight = 60

$ws.Cells.Item(11, 1).Value = "Objectives:"
$ws.Cells.Item(11, 2).Value = $null
$ws.Cells.Item(11, 3).Value = $null
$ws.Rows.Item(11).RowHeight = 60

$ws.Cells.Item(12, 1).Value = "Docentes responsáveis:"
$ws.Cells.Item(12, 2).Value = $null
$ws.Cells.Item(12, 3).Value = $null
$ws.Rows.Item(12).RowHeight = 15

$ws.Cells.Item(13, 1).Value = "Programa resumido:"
$ws.Cells.Item(13, 2).Value = "01/01/2012"
$ws.Cells.Item(13, 3).Value = "01/01/2012"
$ws.Rows.Item(13).RowHeight = 60

$ws.Cells.Item(14, 1).Value = "Short syllabus:"
$ws.Cells.Item(14, 2).Value = $null
$ws.Cells.Item(14, 3).Value = $null
$ws.Rows.Item(14).RowHeight = 60

$ws.Cells.Item(15, 1).Value = "Programa:"
$ws.Cells.Item(15, 2).Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Cells.Item(15, 3).Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Rows.Item(15).RowHeight = 120

$ws.Cells.Item(16, 1).Value = "Syllabus:"
$ws.Cells.Item(16, 2).Value = $null
$ws.Cells.Item(16, 3).Value = $null
$ws.Rows.Item(16).RowHeight = 120

$ws.Cells.Item(17, 1).Value = "Avaliação:"
$ws.Cells.Item(17, 2).Value = $null
$ws.Cells.Item(17, 3).Value = $null
$ws.Rows.Item(17).RowHeight = 15

$ws.Cells.Item(18, 1).Value = "Método:"
$ws.Cells.Item(18, 2).Value = "519033 - Carlos Yujiro Shigue"
$ws.Cells.Item(18, 3).Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(18).RowHeight = 60

$ws.Cells.Item(19, 1).Value = "Critério:"
$ws.Cells.Item(19, 2).Value = "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa."
$ws.Cells.Item(19, 3).Value = "Este curso deverá conter duas avaliações escritas denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa."
$ws.Rows.Item(19).RowHeight = 60

$ws.Cells.Item(20, 1).Value = "Norma de recuperação:"
$ws.Cells.Item(20, 2).Value = "A média do semestre será computada com base na relação: M=(P1+2P2)/3"
$ws.Cells.Item(20, 3).Value = "A média do semestre será computada com base na relação: M=(P1+2P2)/3"
$ws.Rows.Item(20).RowHeight = 60

$ws.Cells.Item(21, 1).Value = "Bibliografia:"
$ws.Cells.Item(21, 2).Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será computada com base na relação: MF=(M+RC)/2"
$ws.Cells.Item(21, 3).Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será computada com base na relação: MF=(M+RC)/2"
$ws.Rows.Item(21).RowHeight = 120
